$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$xlWhole = 1
$xlByRows = 1

$used = $ws.UsedRange
$used.Replace("dnasr281@gmail.com, System", "System, dnasr281@gmail.com", $xlWhole, $xlByRows, $false, $false, $false, $false)
